$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.877.09"
$ws.Range("E2").Value = "  +1.73%  "

$ws.Range("D3").Value = "1.726.82"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9973"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9978"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.20%  "

$ws.Range("E7").Value = "  -0.34%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2591"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.95%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06210"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("D10").Value = "1.728.09"
$ws.Range("E10").Value = "  +0.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.02"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.60%  "

$ws.Range("E12").Value = "  -1.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6078"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.477"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.90%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.02%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9976"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.20%  "

$ws.Range("D17").Value = "26.640.61"
$ws.Range("E17").Value = "  +0.82%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9973"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007184"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.82%  "

$ws.Range("D21").Value = "1.950.99"
$ws.Range("E21").Value = "  +0.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.418"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.45%  "

$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.085"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.51%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.18%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.40%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.774"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.58%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "106.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.58%  "

$ws.Range("E29").Value = "  -1.30%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.949"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.05%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07998"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.41%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.687"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.36%  "

$ws.Range("E33").Value = "  -0.17%  "

$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.598"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.12%  "

$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.010"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.57%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6258"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.04%  "

$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9350"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.69%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.049"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.58%  "

$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.454"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.71%  "

$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9977"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.19%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01499"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.87%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.662"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.12%  "

$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3849"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.09%  "

$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.850"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.21%  "

$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1162"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.51%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05400"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.57%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.906"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.16%  "

$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.26%  "

$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "51.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.48%  "

$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.236"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.33%  "
